{"js": "// 1) \"\u0422\u0440\u0435\u0445\u043e\u0441\u0435\u0432\u043e\u0439 \u043c\u0430\u0433\u043d\u0438\u0442\u043e\u043c\u0435\u0442\u0440 HMC5883L \u043f\u0440\u0438\u043c\u0435\u043d\u044f\u0435\u0442\u0441\u044f \u0434\u043b\u044f \u0438\u0437\u043c\u0435\u0440\u0435\u043d\u0438\u044f \u0441\u043b\u0430\u0431\u044b\u0445 \u043c\u0430\u0433\u043d\u0438\u0442\u043d\u044b\u0445 \u043f\u043e\u043b\u0435\u0439 \u0432 \u0444\u0438\u0437\u0438\u0447\u0435\u0441\u043a\u0438\u0445 \u0438\u0441\u0441\u043b\u0435\u0434\u043e\u0432\u0430\u043d\u0438\u044f\u0445. \"\n//    -> remove the word \"\u0441\u043b\u0430\u0431\u044b\u0445\" (leaving the extra space behind), i.e.\n//    \"...\u0434\u043b\u044f \u0438\u0437\u043c\u0435\u0440\u0435\u043d\u0438\u044f  \u043c\u0430\u0433\u043d\u0438\u0442\u043d\u044b\u0445 \u043f\u043e\u043b\u0435\u0439...\"\nconst body = context.document.body;\n\nconst searchResults = body.search(\n  \"\u0422\u0440\u0435\u0445\u043e\u0441\u0435\u0432\u043e\u0439 \u043c\u0430\u0433\u043d\u0438\u0442\u043e\u043c\u0435\u0442\u0440 HMC5883L \u043f\u0440\u0438\u043c\u0435\u043d\u044f\u0435\u0442\u0441\u044f \u0434\u043b\u044f \u0438\u0437\u043c\u0435\u0440\u0435\u043d\u0438\u044f \u0441\u043b\u0430\u0431\u044b\u0445 \u043c\u0430\u0433\u043d\u0438\u0442\u043d\u044b\u0445 \u043f\u043e\u043b\u0435\u0439 \u0432 \u0444\u0438\u0437\u0438\u0447\u0435\u0441\u043a\u0438\u0445 \u0438\u0441\u0441\u043b\u0435\u0434\u043e\u0432\u0430\u043d\u0438\u044f\u0445.\",\n  { matchCase: true, matchWholeWord: false }\n);\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly one match for the introduction sentence, found \" +\n      searchResults.items.length\n  );\n}\n\nsearchResults.items[0].insertText(\n  \"\u0422\u0440\u0435\u0445\u043e\u0441\u0435\u0432\u043e\u0439 \u043c\u0430\u0433\u043d\u0438\u0442\u043e\u043c\u0435\u0442\u0440 HMC5883L \u043f\u0440\u0438\u043c\u0435\u043d\u044f\u0435\u0442\u0441\u044f \u0434\u043b\u044f \u0438\u0437\u043c\u0435\u0440\u0435\u043d\u0438\u044f  \u043c\u0430\u0433\u043d\u0438\u0442\u043d\u044b\u0445 \u043f\u043e\u043b\u0435\u0439 \u0432 \u0444\u0438\u0437\u0438\u0447\u0435\u0441\u043a\u0438\u0445 \u0438\u0441\u0441\u043b\u0435\u0434\u043e\u0432\u0430\u043d\u0438\u044f\u0445.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// 2) Insert a new paragraph right after the paragraph ending with\n//    \"...\u041a\u0430\u043b\u0438\u0431\u0440\u043e\u0432\u043a\u0430 \u0441\u0432\u043e\u0434\u0438\u0442\u0441\u044f \u043a \u043f\u043e\u0438\u0441\u043a\u0443 \u043c\u0430\u0442\u0440\u0438\u0446\u044b \u043f\u0440\u0435\u043e\u0431\u0440\u0430\u0437\u043e\u0432\u0430\u043d\u0438\u044f \u044d\u0442\u043e\u0433\u043e \u044d\u043b\u043b\u0438\u043f\u0441\u0430 \u0432 \u0441\u0444\u0435\u0440\u0443.\"\nconst targetResults = body.search(\n  \"\u041a\u0430\u043b\u0438\u0431\u0440\u043e\u0432\u043a\u0430 \u0441\u0432\u043e\u0434\u0438\u0442\u0441\u044f \u043a \u043f\u043e\u0438\u0441\u043a\u0443 \u043c\u0430\u0442\u0440\u0438\u0446\u044b \u043f\u0440\u0435\u043e\u0431\u0440\u0430\u0437\u043e\u0432\u0430\u043d\u0438\u044f \u044d\u0442\u043e\u0433\u043e \u044d\u043b\u043b\u0438\u043f\u0441\u0430 \u0432 \u0441\u0444\u0435\u0440\u0443.\",\n  { matchCase: true, matchWholeWord: false }\n);\ntargetResults.load(\"text,paragraphs\");\nawait context.sync();\n\nif (targetResults.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly one match for the calibration sentence, found \" +\n      targetResults.items.length\n  );\n}\n\nconst targetParagraph = targetResults.items[0].paragraphs.getFirst();\ntargetParagraph.insertParagraph(\n  \"\u0410\u043b\u0433\u043e\u0440\u0438\u0442\u043c \u0442\u0430\u043a\u0436\u0435 \u043f\u043e\u0434\u0445\u043e\u0434\u0438\u0442 \u0434\u043b\u044f \u043a\u0430\u043b\u0438\u0431\u0440\u043e\u0432\u043a\u0438 \u0434\u0440\u0443\u0433\u0438\u0445 \u0442\u0440\u0435\u0445\u043e\u0441\u0435\u0432\u044b\u0445 \u043c\u0430\u0433\u043d\u0438\u0442\u043e\u043c\u0435\u0442\u0440\u043e\u0432, \u0430 \u0442\u0430\u043a\u0436\u0435 \u0434\u043b\u044f \u043f\u0440\u043e\u0432\u0435\u0440\u043a\u0438 \u0442\u043e\u0447\u043d\u043e\u0441\u0442\u0438 \u043a\u0430\u043b\u0438\u0431\u0440\u043e\u0432\u043a\u0438 \u0442\u0435\u0445, \u043a\u043e\u0442\u043e\u0440\u044b\u0435 \u043f\u043e\u0441\u0442\u0430\u0432\u043b\u044f\u044e\u0442\u0441\u044f \u0432 \u043e\u0442\u043a\u0430\u043b\u0438\u0431\u0440\u043e\u0432\u0430\u043d\u043d\u043e\u043c \u0432\u0438\u0434\u0435.\",\n  Word.InsertLocation.after\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------------\n# 1) \"...\u043f\u0440\u0438\u043c\u0435\u043d\u044f\u0435\u0442\u0441\u044f \u0434\u043b\u044f \u0438\u0437\u043c\u0435\u0440\u0435\u043d\u0438\u044f \u0441\u043b\u0430\u0431\u044b\u0445 \u043c\u0430\u0433\u043d\u0438\u0442\u043d\u044b\u0445 \u043f\u043e\u043b\u0435\u0439 \u0432 \u0444\u0438\u0437\u0438\u0447\u0435\u0441\u043a\u0438\u0445\n#    \u0438\u0441\u0441\u043b\u0435\u0434\u043e\u0432\u0430\u043d\u0438\u044f\u0445. \" -> drop the word \"\u0441\u043b\u0430\u0431\u044b\u0445\" (the space on each side of\n#    it stays, producing a double space between \"\u0438\u0437\u043c\u0435\u0440\u0435\u043d\u0438\u044f\" and \"\u043c\u0430\u0433\u043d\u0438\u0442\u043d\u044b\u0445\").\n# ---------------------------------------------------------------------------\n$find1 = $d.Content.Find\n$find1.ClearFormatting()\n$find1.Text = \"\u0441\u043b\u0430\u0431\u044b\u0445\"\n$found1 = $find1.Execute()\nif (-not $found1) {\n  throw \"Could not find the word '\u0441\u043b\u0430\u0431\u044b\u0445' in the introduction paragraph.\"\n}\n$find1.Parent.Text = \"\"\n\n# ---------------------------------------------------------------------------\n# 2) Insert a brand-new paragraph (same Style17 style as its neighbours)\n#    right after the paragraph that ends with \"...\u041a\u0430\u043b\u0438\u0431\u0440\u043e\u0432\u043a\u0430 \u0441\u0432\u043e\u0434\u0438\u0442\u0441\u044f \u043a\n#    \u043f\u043e\u0438\u0441\u043a\u0443 \u043c\u0430\u0442\u0440\u0438\u0446\u044b \u043f\u0440\u0435\u043e\u0431\u0440\u0430\u0437\u043e\u0432\u0430\u043d\u0438\u044f \u044d\u0442\u043e\u0433\u043e \u044d\u043b\u043b\u0438\u043f\u0441\u0430 \u0432 \u0441\u0444\u0435\u0440\u0443.\"\n# ---------------------------------------------------------------------------\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"\u041a\u0430\u043b\u0438\u0431\u0440\u043e\u0432\u043a\u0430 \u0441\u0432\u043e\u0434\u0438\u0442\u0441\u044f \u043a \u043f\u043e\u0438\u0441\u043a\u0443 \u043c\u0430\u0442\u0440\u0438\u0446\u044b \u043f\u0440\u0435\u043e\u0431\u0440\u0430\u0437\u043e\u0432\u0430\u043d\u0438\u044f \u044d\u0442\u043e\u0433\u043e \u044d\u043b\u043b\u0438\u043f\u0441\u0430 \u0432 \u0441\u0444\u0435\u0440\u0443.\"\n$found2 = $find2.Execute()\nif (-not $found2) {\n  throw \"Could not find the calibration-summary sentence.\"\n}\n$targetRange = $find2.Parent\n\n# Locate the 1-based Paragraphs index of the paragraph holding that sentence.\n$targetParaIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n  $p = $d.Paragraphs.Item($i).Range\n  if ($p.Start -le $targetRange.Start -and $p.End -ge $targetRange.End) {\n    $targetParaIndex = $i\n    break\n  }\n}\nif ($targetParaIndex -eq -1) {\n  throw \"Could not resolve the paragraph index for the calibration-summary sentence.\"\n}\n\n$targetRange.Collapse(0)\n$targetRange.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Item($targetParaIndex + 1)\n$newPara.Range.InsertBefore(\"\u0410\u043b\u0433\u043e\u0440\u0438\u0442\u043c \u0442\u0430\u043a\u0436\u0435 \u043f\u043e\u0434\u0445\u043e\u0434\u0438\u0442 \u0434\u043b\u044f \u043a\u0430\u043b\u0438\u0431\u0440\u043e\u0432\u043a\u0438 \u0434\u0440\u0443\u0433\u0438\u0445 \u0442\u0440\u0435\u0445\u043e\u0441\u0435\u0432\u044b\u0445 \u043c\u0430\u0433\u043d\u0438\u0442\u043e\u043c\u0435\u0442\u0440\u043e\u0432, \u0430 \u0442\u0430\u043a\u0436\u0435 \u0434\u043b\u044f \u043f\u0440\u043e\u0432\u0435\u0440\u043a\u0438 \u0442\u043e\u0447\u043d\u043e\u0441\u0442\u0438 \u043a\u0430\u043b\u0438\u0431\u0440\u043e\u0432\u043a\u0438 \u0442\u0435\u0445, \u043a\u043e\u0442\u043e\u0440\u044b\u0435 \u043f\u043e\u0441\u0442\u0430\u0432\u043b\u044f\u044e\u0442\u0441\u044f \u0432 \u043e\u0442\u043a\u0430\u043b\u0438\u0431\u0440\u043e\u0432\u0430\u043d\u043d\u043e\u043c \u0432\u0438\u0434\u0435.\")\n"}
